$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = "Jordan Poole"
$ws.Cells.Item(2, 2).Value = "PG,SG"
$ws.Cells.Item(2, 3).Value = "Washington Wizards"
$ws.Cells.Item(3, 1).Value = "RJ Barrett"
$ws.Cells.Item(3, 2).Value = "SG,SF,PF"
$ws.Cells.Item(3, 3).Value = "Toronto Raptors"
$ws.Cells.Item(4, 1).Value = "Tobias Harris"
$ws.Cells.Item(4, 2).Value = "SF,PF"
$ws.Cells.Item(4, 3).Value = "Detroit Pistons"
$ws.Cells.Item(5, 1).Value = "Nikola Jovic"
$ws.Cells.Item(5, 2).Value = "PF,C"
$ws.Cells.Item(5, 3).Value = "Miami Heat"
$ws.Cells.Item(6, 1).Value = "Jose Alvarado"
$ws.Cells.Item(6, 2).Value = "PG"
$ws.Cells.Item(6, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(7, 1).Value = "CJ McCollum"
$ws.Cells.Item(7, 2).Value = "PG,SG"
$ws.Cells.Item(7, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(8, 1).Value = "Jalen Williams"
$ws.Cells.Item(8, 2).Value = "SG,SF,PF,C"
$ws.Cells.Item(8, 3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(9, 1).Value = "Guerschon Yabusele"
$ws.Cells.Item(9, 2).Value = "PF,C"
$ws.Cells.Item(9, 3).Value = "Philadelphia 76ers"
$ws.Cells.Item(10, 1).Value = "Quentin Grimes"
$ws.Cells.Item(10, 2).Value = "SG,SF"
$ws.Cells.Item(10, 3).Value = "Dallas Mavericks"
$ws.Cells.Item(11, 1).Value = "Zach LaVine"
$ws.Cells.Item(11, 2).Value = "SG,SF"
$ws.Cells.Item(11, 3).Value = "Chicago Bulls"
$ws.Cells.Item(12, 1).Value = "Lauri Markkanen"
$ws.Cells.Item(12, 2).Value = "SF,PF"
$ws.Cells.Item(12, 3).Value = "Utah Jazz"
$ws.Cells.Item(13, 1).Value = "Shai Gilgeous-Alexander"
$ws.Cells.Item(13, 2).Value = "PG,SG"
$ws.Cells.Item(13, 3).Value = "Oklahoma City Thunder"
$ws.Cells.Item(14, 1).Value = "Andrew Nembhard"
$ws.Cells.Item(14, 2).Value = "PG,SG"
$ws.Cells.Item(14, 3).Value = "Indiana Pacers"
$ws.Cells.Item(15, 1).Value = "Joel Embiid"
$ws.Cells.Item(15, 2).Value = "C"
$ws.Cells.Item(15, 3).Value = "Philadelphia 76ers"
$ws.Cells.Item(16, 1).Value = "Toumani Camara"
$ws.Cells.Item(16, 2).Value = "SF,PF"
$ws.Cells.Item(16, 3).Value = "Portland Trail Blazers"
$ws.Cells.Item(17, 1).Value = "John Collins"
$ws.Cells.Item(17, 2).Value = "PF,C"
$ws.Cells.Item(17, 3).Value = "Utah Jazz"
$ws.Cells.Item(18, 1).Value = "Kyrie Irving"
$ws.Cells.Item(18, 2).Value = "PG,SG"
$ws.Cells.Item(18, 3).Value = "Dallas Mavericks"
$ws.Cells.Item(19, 1).Value = "Jimmy Butler"
$ws.Cells.Item(19, 2).Value = "SF,PF"
$ws.Cells.Item(19, 3).Value = "Miami Heat"
